# Add a new worksheet "SurveyRespnod" capturing survey/responder timings,
# following the same general layout used on the other timing sheets
# (e.g. "PushPull"), per the commit:
#   "Make a new spreadsheet page on which to capture timings."

$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook (after the last existing sheet).
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "SurveyRespnod"

# Row 1: title cells
$ws.Range("A1").Value = "tcp:"
$ws.Range("B1").Value = "10000 survyeys"

# Row 3: column headers (written before row 2 so the shared-string table
# picks up "Size"/"big surveys/sec"/etc. ahead of the "N responders" labels)
$ws.Range("A3").Value = "Size"
$ws.Range("B3").Value = "big surveys/sec"
$ws.Range("C3").Value = "big KBs"
$ws.Range("D3").Value = "small surveys/sec"
$ws.Range("E3").Value = "small KBs"
$ws.Range("F3").Value = "big surveys/sec"
$ws.Range("G3").Value = "big KBs"
$ws.Range("H3").Value = "small surveys/sec"
$ws.Range("I3").Value = "small KBs"
$ws.Range("J3").Value = "big surveys/sec"
$ws.Range("K3").Value = "big KBs"
$ws.Range("L3").Value = "small surveys/sec"
$ws.Range("M3").Value = "small KBs"
$ws.Range("N3").Value = "small surveys/sec"
$ws.Range("O3").Value = "small KBs"
$ws.Range("P3").Value = "big surveys/sec"
$ws.Range("Q3").Value = "big KBs"
$ws.Range("R3").Value = "small surveys/sec"
$ws.Range("S3").Value = "small KBs"

# Row 2: responder-count group headers
$ws.Range("C2").Value = "1 responder"
$ws.Range("G2").Value = "2 responders"
$ws.Range("K2").Value = "3 responders"
$ws.Range("O2").Value = "4responders"

# Row 4: sample timing data for the smallest message size
$ws.Range("A4").Value = 1024
$ws.Range("B4").Value = 1935
$ws.Range("C4").Value = 1935
$ws.Range("D4").Value = 2284
$ws.Range("E4").Value = 2284
$ws.Range("F4").Value = 1250
$ws.Range("G4").Value = 1250
$ws.Range("H4").Value = 1371
$ws.Range("I4").Value = 2742

# Column A: message size progression (plain values then doubling formulas)
$ws.Range("A5").Value = 2048
$ws.Range("A6").Value = 4096
$ws.Range("A7").Value = 16384
$ws.Range("A8").Value = 32768
$ws.Range("A9").Formula = "=A8*2"
$ws.Range("A10").Formula = "=A9*2"
$ws.Range("A11").Formula = "=A10*2"
$ws.Range("A12").Formula = "=A11*2"
$ws.Range("A13").Formula = "=A12*2"

# Approximate the "best fit" column widths Excel computed for the header text.
$ws.Columns.Item(2).ColumnWidth = 13.6415
$ws.Columns.Item(4).ColumnWidth = 15.9835
$ws.Columns.Item(5).ColumnWidth = 8.833
$ws.Columns.Item(6).ColumnWidth = 13.6415
$ws.Columns.Item(7).ColumnWidth = 11.6415
$ws.Columns.Item(8).ColumnWidth = 15.9835
$ws.Columns.Item(9).ColumnWidth = 8.833

# Match the selection left active on the new sheet and make it the active tab.
$ws.Range("N3:S3").Select() | Out-Null
$ws.Activate()
